$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F79").Value = 69
$ws.Range("G79").Value = 4289.73
$ws.Range("B90").Value = 160978.74
$ws.Range("B112").Value = 57756
$ws.Range("E112").Value = 79.37
$ws.Range("F112").Value = -100
$ws.Range("G112").Value = -6644
$ws.Range("B113").Value = 64350
$ws.Range("E113").Value = 70.63
$ws.Range("F113").Value = 2
$ws.Range("G113").Value = 132.88
$ws.Range("B127").Value = 64329
$ws.Range("E127").Value = 128.32
$ws.Range("F127").Value = 1
$ws.Range("G127").Value = 120.69
$ws.Range("B128").Value = 57552
$ws.Range("E128").Value = 136.86
$ws.Range("F128").Value = -5
$ws.Range("G128").Value = -603.45
$ws.Range("F144").Value = 909
$ws.Range("G144").Value = 7681.05
$ws.Range("B147").Value = 11800.03
$ws.Range("F163").Value = 10
$ws.Range("G163").Value = 2659.6
$ws.Range("B175").Value = 25955.49
$ws.Range("F205").Value = 14
$ws.Range("G205").Value = 5279.96
$ws.Range("B216").Value = 29333.52
$ws.Range("B227").Value = 55373
$ws.Range("E227").Value = 163.62
$ws.Range("F227").Value = -94
$ws.Range("G227").Value = -13562.32
$ws.Range("B228").Value = 63520
$ws.Range("E228").Value = 153.4
$ws.Range("F228").Value = 65
$ws.Range("G228").Value = 9378.2
$ws.Range("B232").Value = 63510
$ws.Range("E232").Value = 50.66
$ws.Range("F232").Value = 112
$ws.Range("G232").Value = 5335.68
$ws.Range("B233").Value = 55356
$ws.Range("E233").Value = 54.04
$ws.Range("F233").Value = -158
$ws.Range("G233").Value = -7527.12
$ws.Range("B243").Value = 63560
$ws.Range("E243").Value = 134.87
$ws.Range("F243").Value = 1
$ws.Range("G243").Value = 126.86
$ws.Range("B244").Value = 60325
$ws.Range("E244").Value = 151.57
$ws.Range("F244").Value = -102
$ws.Range("G244").Value = -12939.72
$ws.Range("F338").Value = 68
$ws.Range("G338").Value = 1611.6
$ws.Range("B346").Value = 22874.07
$ws.Range("B382").Value = 65066
$ws.Range("E382").Value = 13.61
$ws.Range("F382").Value = 90
$ws.Range("G382").Value = 1152.9
$ws.Range("B383").Value = 53263
$ws.Range("E383").Value = 15.29
$ws.Range("F383").Value = -309
$ws.Range("G383").Value = -3958.29
$ws.Range("B391").Value = 45718
$ws.Range("E391").Value = 19.38
$ws.Range("F391").Value = -294
$ws.Range("G391").Value = -4768.68
$ws.Range("B392").Value = 64927
$ws.Range("E392").Value = 17.26
$ws.Range("F392").Value = 106
$ws.Range("G392").Value = 1719.32
$ws.Range("B398").Value = 45702
$ws.Range("E398").Value = 31.43
$ws.Range("F398").Value = -215
$ws.Range("G398").Value = -5654.5
$ws.Range("B399").Value = 64919
$ws.Range("E399").Value = 27.97
$ws.Range("F399").Value = 61
$ws.Range("G399").Value = 1604.3
$ws.Range("B401").Value = 65067
$ws.Range("E401").Value = 15.65
$ws.Range("F401").Value = 126
$ws.Range("G401").Value = 1855.98
$ws.Range("B402").Value = 53595
$ws.Range("E402").Value = 17.61
$ws.Range("F402").Value = -335
$ws.Range("G402").Value = -4934.55
$ws.Range("B458").Value = 53319
$ws.Range("E458").Value = 310.64
$ws.Range("F458").Value = -6
$ws.Range("G458").Value = -1643.52
$ws.Range("B459").Value = 64810
$ws.Range("E459").Value = 291.22
$ws.Range("F459").Value = 4
$ws.Range("G459").Value = 1095.68
$ws.Range("B489").Value = 64830
$ws.Range("E489").Value = 34.9
$ws.Range("F489").Value = 104
$ws.Range("G489").Value = 3414.32
$ws.Range("B490").Value = 60022
$ws.Range("E490").Value = 37.22
$ws.Range("F490").Value = -113
$ws.Range("G490").Value = -3709.79
$ws.Range("B594").Value = 65079
$ws.Range("F594").Value = 6
$ws.Range("G594").Value = 245.22
$ws.Range("B595").Value = 65362
$ws.Range("F595").Value = 18
$ws.Range("G595").Value = 735.66
$ws.Range("F600").Value = 51
$ws.Range("G600").Value = 2544.39
$ws.Range("B605").Value = 11619.99
$ws.Range("B647").Value = 2264685.95
$ws.Range("B648").Value = 2264685.95
